$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape "Rectangle 23" (id=24): "OPNFV test tooling / RC-2 test suite"
#     -> "Anuket" + " test tooling / RC-2 test suite"
$sh1 = $s.Shapes.Item(1)
$tr1 = $sh1.TextFrame.TextRange
$c1 = $tr1.Characters(1, 5)          # "OPNFV"
$c1.Text = "Anuket"

$tr1b = $sh1.TextFrame.TextRange
$newRun1 = $tr1b.Characters(1, 6)    # "Anuket"
$newRun1.LanguageID = "LID4096"

# --- Shape "Rectangle 13" (id=14): "OPNFV Kuberef: RI-2 integration"
#     -> "Anuket" + " " + "Kuberef" + ": RI-2 integration"
$sh2 = $s.Shapes.Item(3)
$tr2 = $sh2.TextFrame.TextRange
$c2 = $tr2.Characters(1, 5)          # "OPNFV"
$c2.Text = "Anuket"

$tr2b = $sh2.TextFrame.TextRange
$newRun2 = $tr2b.Characters(1, 6)    # "Anuket"
$newRun2.LanguageID = "LID4096"
